# Generate Report for Handoff
# - Update status text from "Handed back: in sync with en-US" to "Ready for handoff"
#   on the Overview sheet (zh-cn/de-de status columns) and on each language sheet's
#   "Status" column.
# - Refresh the handoff/generate timestamps to reflect the new handoff run.
# - Narrow the "Latest HO Xliff Generate Date" / "Status" columns that used to be
#   sized for the old, longer status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Cell value updates -----------------------------------------------------

# Overview sheet: zh-cn (E2) / de-de (F2) status, and the generate date (G2)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-20 05:01:38"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (H2)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-20 05:01:34"

# de-de sheet: Status (C2) and Latest Handoff Datetime (H2)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-20 05:01:38"

# --- Column width updates ---------------------------------------------------
# Old width (29.9777047293527 chars) shrinks to 17.2159881591797 chars now that
# the status column no longer needs to fit the long "Handed back..." text.
# The host's ColumnWidth setter snaps to whole-pixel character units, so use
# the closest representable width (16.3333... => stored width ~17.16667).
$newColWidth = 16.333333333333332

$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth   # Overview!E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth   # Overview!F (de-de status)
$wsZhCn.Columns.Item(3).ColumnWidth     = $newColWidth   # zh-cn!C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth     = $newColWidth   # de-de!C (Status)
